# Regenerate merged AHB files
# - Rename the diff header columns from the "_old"/"_new" comparison
#   labels to the concrete version labels "_FV2310"/"_FV2404".
# - Freeze the header row.
# - Turn the used range into a proper Excel Table ("Table1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (columns A:J are the "old"/FV2310 side,
#        L:U are the "new"/FV2404 side; K holds the "diff" column and is
#        left untouched). ---
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"

$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# --- 2. Freeze the header row (row 1) ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null

# --- 3. Turn A1:U72 into an Excel Table ("Table1") with a plain style
#        (no extra header shading beyond the sheet's own header format). ---
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U72"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = ""
